$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing header row, pushing "NIM" down to row 3
$ws.Rows("1:1").Insert()
$ws.Rows("1:1").Insert()

# New "Group Name " header in A1, spanning visually across A1:B1
$ws.Range("A1").Value = "Group Name "
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Font.Size = 16

# Copy A1's formatting (bold, 16pt font) onto B1 so both cells share style s=2
$a1.Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights: title row taller, spacer/NIM rows standard
$ws.Rows(1).RowHeight = 21
$ws.Rows(2).RowHeight = 15

# Widen columns A and B
$ws.Range("A1:B1").ColumnWidth = 20.77734375
